# Update "Overall" weighted evaluation table (rows 3-10) on Sayfa1 with the
# actual per-reviewer scores. Columns: A=raw weight input, D/E/F/G=criteria
# scores for the four evaluators. B (fraction) and downstream
# tables recalc automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

# Row 3 - Smart Connected Cat Feeding & Monitoring System ("Gimme Fast")
$ws.Range("A3").Value = 3
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 10

# Row 4 - Autonomous Valet Parking Service
$ws.Range("A4").Value = 7
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2

# Row 5 - Gimme Fast
$ws.Range("A5").Value = 6
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2

# Row 6 - Where am I
$ws.Range("A6").Value = 8
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3

# Row 7
$ws.Range("A7").Value = 7
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4

# Row 8
$ws.Range("A8").Value = 5
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 9

# Row 9
$ws.Range("A9").Value = 4
$ws.Range("D9").Value = 6
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3

# Row 10
$ws.Range("A10").Value = 10
$ws.Range("D10").Value = 10
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 7

$excel.Calculate()
